$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "62×91=" "19×85="
Replace-Text "98×53=" "77×64="
Replace-Text "83×62=" "53×91="
Replace-Text "35×73=" "62×44="
Replace-Text "42×44=" "85×58="
Replace-Text "41×12=" "14×40="
Replace-Text "26×19=" "22×65="
Replace-Text "91×63=" "70×53="
Replace-Text "40×96=" "73×80="
Replace-Text "21×96=" "76×64="
Replace-Text "50×25=" "12×96="
Replace-Text "46×36=" "36×91="
Replace-Text "83×19=" "21×87="
Replace-Text "23×62=" "51×34="
Replace-Text "61×87=" "16×96="
Replace-Text "89×60=" "31×74="
Replace-Text "51×92=" "73×80="
Replace-Text "18×99=" "89×19="
Replace-Text "77×99=" "94×79="
Replace-Text "85×80=" "42×51="
Replace-Text "40×76=" "90×21="
Replace-Text "44×65=" "78×79="
Replace-Text "16×27=" "14×26="
Replace-Text "65×88=" "26×53="
Replace-Text "82×18=" "90×24="
